$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

# Copy formatting (date style) from the cell above so the new date cell
# matches the existing column A formatting (style index 1 / numFmtId 22).
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 42622.888495370367
$ws.Cells.Item($row, 2).Value = -8
$ws.Cells.Item($row, 3).Value = 55
$ws.Cells.Item($row, 4).Value = 42
$ws.Cells.Item($row, 5).Value = 55
$ws.Cells.Item($row, 6).Value = 24
$ws.Cells.Item($row, 7).Value = 19681
$ws.Cells.Item($row, 8).Value = 13408
$ws.Cells.Item($row, 9).Value = 731
$ws.Cells.Item($row, 10).Value = 159
$ws.Cells.Item($row, 11).Value = 121
$ws.Cells.Item($row, 12).Value = 18
$ws.Cells.Item($row, 13).Value = 6
$ws.Cells.Item($row, 14).Value = "Named"

$wb.Save()
